$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

# Sweet 16 "First Four" date column shifts up a day across four matchups
# (the missing "Powdered Mini Donuts" slot collapsed the schedule by a day).
$ws.Range("D8").Value = 44663
$ws.Range("D13").Value = 44664
$ws.Range("D18").Value = 44665
$ws.Range("D23").Value = 44666

# The two "Powdered Mini Donuts" / "Chocolate Mini Donuts" seeds (15/16) swap
# places in the bracket, as do the "Peanut Butter Crunch" / "Chocolate Chip
# Cream Pies" seeds (17/18). The XLOOKUP formulas in column B recalc on their own.
$ws.Range("A9").Value = 16
$ws.Range("A12").Value = 17
$ws.Range("A29").Value = 15
$ws.Range("A32").Value = 18

# First Four matchup date for the Chocolate Chip Cream Pies / TBD slot moves out.
$ws.Range("B30").Value = 44669
